$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New image links replacing the old ones in column B (rows 2-12)
$ws.Range("B2").ClearFormats()
$ws.Range("B2").Value = "https://imgbb.host/images/QP4Fq.png"

$ws.Range("B3").ClearFormats()
$ws.Range("B3").Value = "https://imgbb.host/images/QPWMS.png"

$ws.Range("B4").ClearFormats()
$ws.Range("B4").Value = "https://imgbb.host/images/QPyE7.png"

$ws.Range("B5").ClearFormats()
$ws.Range("B5").Value = "https://imgbb.host/images/QPAru.png"

$ws.Range("B6").ClearFormats()
$ws.Range("B6").Value = "https://imgbb.host/images/QPxZV.png"

$ws.Range("B7").ClearFormats()
$ws.Range("B7").Value = "https://imgbb.host/images/QPdzb.png"

$ws.Range("B8").ClearFormats()
$ws.Range("B8").Value = "https://imgbb.host/images/QPIhB.png"

$ws.Range("B9").ClearFormats()
$ws.Range("B9").Value = "https://imgbb.host/images/QPiUi.png"

$ws.Range("B10").ClearFormats()
$ws.Range("B10").Value = "https://imgbb.host/images/QPE0P.png"

$ws.Range("B11").Value = "https://imgbb.host/images/QPP4E.png"

$ws.Range("B12").ClearFormats()
$ws.Range("B12").Value = "https://imgbb.host/images/QPXmM.png"

# Update the active selection to B13
$ws.Range("B13").Select()
